$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - update "interested count" (F column) figures,
# and flip row 5's min-price (G) from the inlineStr placeholder "不可售"
# (not for sale) to an actual numeric price now that tickets are on sale.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$expoF = @{
    2  = 216
    3  = 54747
    4  = 1341
    5  = 379
    6  = 320
    7  = 880
    8  = 756
    9  = 398
    10 = 3058
    11 = 904
    12 = 5228
    13 = 1283
    14 = 1015
    18 = 403
    19 = 1280
    20 = 98
    22 = 177
    23 = 363
    24 = 25
    28 = 61
    29 = 5035
    31 = 4960
    32 = 8959
    35 = 135
    36 = 221
    37 = 426
    38 = 114
    40 = 4205
    41 = 244
}

foreach ($row in $expoF.Keys) {
    $wsExpo.Range("F$row").Value = $expoF[$row]
}

# Row 5 ticket now has a real minimum price instead of "not for sale"
$wsExpo.Range("G5").Value = 58

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F12").Value = 1132

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 138
$wsLocal.Range("F5").Value = 40

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - a brand new event ("逐月节") was inserted
# at row 4, pushing the following four events down by one row each; the
# event that used to sit at row 8 dropped off this sheet entirely, and every
# row from 9 onward keeps its own identity with just an updated F figure.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F3").Value = 1341

# New row 4: previously-unlisted event, now inserted at the top of the
# 10-04/10-05 cluster.
# (B4 holds a date formatted as *text*, e.g. "2024-10-04"; without forcing a
# text number format first, Excel's smart-entry would silently reinterpret
# it as a date serial number. Reset back to Normal style afterwards so we
# don't leave a stray number-format override on the cell.)
$wsAll.Range("B4").NumberFormat = "@"
$wsAll.Range("B4").Value = "2024-10-04"
$wsAll.Range("B4").Style = "Normal"
$wsAll.Range("C4").Value = "杭州·逐月节·园游会·原神×绝区零×崩铁×崩坏同人only"
$wsAll.Range("D4").Value = "莫干山路987号 资辉壹方汇"
$wsAll.Range("E4").Value = "2024.10.04 09:30-10.05 17:00"
$wsAll.Range("F4").Value = 379
$wsAll.Range("G4").Value = 58
$wsAll.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=92406"
$wsAll.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202409/mQh43oPd1726134932363.png"

# Row 5 now holds what used to be row 4's event.
$wsAll.Range("C5").Value = "杭州·2024·华彩的摔跤宴 - 木吉KAZUYA降临"
$wsAll.Range("D5").Value = "莫干山路188-200号 之江饭店(莫干山路店)"
$wsAll.Range("E5").Value = "2024.10.05 10:00-10.05 16:00"
$wsAll.Range("F5").Value = 320
$wsAll.Range("G5").Value = 128
$wsAll.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=92402"
$wsAll.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202409/ZylQGk1P1726033298213.png"

# Row 6 now holds what used to be row 5's event.
$wsAll.Range("C6").Value = "杭州·文豪野犬同人only2.0"
$wsAll.Range("D6").Value = "望江东路333号 杭州瑞莱克斯大酒店"
$wsAll.Range("E6").Value = "2024.10.05 10:00-10.05 17:00"
$wsAll.Range("F6").Value = 880
$wsAll.Range("G6").Value = 54
$wsAll.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=92226"
$wsAll.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202409/nkCZvaiO1725872765608.jpeg"

# Row 7 now holds what used to be row 6's event.
$wsAll.Range("C7").Value = "杭州·次元幻想【玩美大舞台&吃/换谷大会】（免费活动）"
$wsAll.Range("D7").Value = "文三路 玩美的一天沉浸式生活街区"
$wsAll.Range("E7").Value = "2024.10.05 10:00-10.05 17:00"
$wsAll.Range("F7").Value = 756
$wsAll.Range("G7").Value = 30
$wsAll.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=92028"
$wsAll.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202409/FaEB96xH1725394323651.jpeg"

# Row 8 now holds what used to be row 7's event (row 7's former content,
# the "燃梦星辰" event, drops off the sheet).
$wsAll.Range("C8").Value = "杭州·火影同人ONLY"
$wsAll.Range("D8").Value = "5号大街297号 盛泰开元名都大酒店"
$wsAll.Range("E8").Value = "2024.10.05 10:00-10.05 18:00"
$wsAll.Range("F8").Value = 398
$wsAll.Range("G8").Value = 78
$wsAll.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=92458"
$wsAll.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202409/5choDLVP1726713753891.png"

# Remaining rows keep their own event identity; only the F ("want to go"
# count) figure is refreshed.
$allF = @{
    9  = 904
    11 = 1283
    12 = 40
    14 = 1015
    17 = 403
    19 = 1280
    21 = 98
    22 = 177
    24 = 363
    25 = 25
    27 = 61
    28 = 5036
    30 = 8959
    34 = 135
    35 = 221
    36 = 426
    39 = 114
    41 = 4205
    48 = 244
}

foreach ($row in $allF.Keys) {
    $wsAll.Range("F$row").Value = $allF[$row]
}
